# "moved from static methods to classes" -> fill in the USERINFO sheet's
# login/password cells: B2 gets the email login (as a live mailto
# hyperlink, keeping its pre-existing "Hyperlink" cell style), B3 gets
# the password, and the active selection moves to B3.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("USERINFO")

# Stash B2's current formatting (the "Hyperlink" style it already carries)
# in a scratch cell so we can restore it after Hyperlinks.Add, which
# otherwise mints a brand-new (duplicate) cell style.
$ws.Range("D1").Value = "__fmt_scratch__"
$ws.Range("B2").Copy($ws.Range("D1"))

$ws.Range("B2").Value = "alex-borrow@mail.ru"
$ws.Hyperlinks.Add($ws.Range("B2"), "mailto:alex-borrow@mail.ru")

$ws.Range("D1").Copy()
$ws.Range("B2").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("D1").Clear()

$ws.Range("B3").Value = "eGeEVSckqkVGee8VwWvc"

$ws.Range("B3").Select() | Out-Null
